# Update the "17"/"22" id suffixes used to build test-data usernames/emails
# to "18"/"23" respectively (mirrors the source commit which bumped the
# generated-data batch numbers), across all three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sheet1": the two numeric seed cells drive the CONCATENATE()
# formulas for the rest of the block (columns A/B/C), so updating them
# is all that's needed there - Excel recalculates the cached formula
# results automatically.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Range("I2").Value = 18
$wsData.Range("I23").Value = 23

# Restore the scrolled viewport on "Sheet1" (best effort - selection
# stays as-is, only the top-left visible cell moves down a few rows).
# The previously-active sheet/tab is restored afterwards so this is a
# no-op with respect to which sheet/tab is active in the saved file.
$originalActiveSheet = $wb.ActiveSheet
$wsData.Activate()
$wsData.Range("A23:C42").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
if ($originalActiveSheet.Name -ne $wsData.Name) {
    $originalActiveSheet.Activate()
}

# ---------------------------------------------------------------------
# Sheet "login": columns G/H hold the plain name, column I holds the
# name + "@gmail.com". These are static (non-formula) cached values
# that mirror the first 20-row block from "Sheet1", so they need to be
# rewritten explicitly.
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("login")
$loginNames = @(
    "EthanBaker","DelanieCarman","BretAgnew","EdgardoTaylor","TyrekReis",
    "LeannaChow","TuckerCarlson","AnnmarieConnor","MoniqueWitte","MikelWhitlock",
    "VincentAmaya","KeiraQuiroz","EllisCreech","DionteCreel","NicholeFoust",
    "ManuelConnell","LourdesElam","LincolnFrederick","AlisaCash","LucilleGriffiths"
)
for ($i = 0; $i -lt $loginNames.Length; $i++) {
    $row = 2 + $i
    $newName = "$($loginNames[$i])18"
    $wsLogin.Range("G$row").Value = $newName
    $wsLogin.Range("H$row").Value = $newName
    $wsLogin.Range("I$row").Value = "$newName@gmail.com"
}

# ---------------------------------------------------------------------
# Sheet "order": columns R/S hold the plain name, column T holds the
# name + "@gmail.com", mirroring the second 20-row block from "Sheet1".
# ---------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item("order")
$orderNames = @(
    "DonnellJernigan","MalikOtoole","AlanCaudill","AdanApplegate","AiyanaWhitworth",
    "MercedezBrien","DuaneHager","LorenBell","GeraldHiller","DeionBranch",
    "DakotaHalstead","ElliottFurman","MiltonCamp","DawnChester","ZacheryPetrie",
    "EstebanAngel","JimmyBlankenship","AllysaGrice","AugustineYoo","BrandiSouthard"
)
for ($i = 0; $i -lt $orderNames.Length; $i++) {
    $row = 2 + $i
    $newName = "$($orderNames[$i])23"
    $wsOrder.Range("R$row").Value = $newName
    $wsOrder.Range("S$row").Value = $newName
    $wsOrder.Range("T$row").Value = "$newName@gmail.com"
}
